# Sync evaluation-by-extraction.py and evaluation-by-generation.py: exception handling!
#
# The evaluation table gains a new leading "type" column (all rows tagged
# "generation"), the existing metric columns are renamed, and four new
# trailing columns are added to capture exception/error counts produced by
# the updated exception handling in the python scripts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing A:E data one column to the right (B:F), preserving
# values/styles, to make room for the new leading "type" column.
$ws.Columns.Item(1).Insert()

# --- Header row -----------------------------------------------------------
$ws.Range("A1").Value = "type"
$ws.Range("B1").Value = "model"
$ws.Range("C1").Value = "prec_mean"
$ws.Range("D1").Value = "rec_mean"
$ws.Range("E1").Value = "f1_mean"
$ws.Range("F1").Value = "valid_count"
$ws.Range("G1").Value = "invalid_count"
$ws.Range("H1").Value = "JSON range error"
$ws.Range("I1").Value = "JSON format error"
$ws.Range("J1").Value = "JSON key error"

# New header cells (A1 lost its format when the column was inserted; G1:J1
# are brand new) should look like the rest of the bold/boxed header row.
# (Applied as two separate Range objects -- a single multi-area Range only
# picks up formatting on its first area here.)
$newHeaderCellA = $ws.Range("A1")
$newHeaderCellA.Font.Bold = $true
$newHeaderCellA.HorizontalAlignment = -4108  # xlCenter
$newHeaderCellA.VerticalAlignment = -4160    # xlTop
$newHeaderCellA.Borders.LineStyle = 1        # xlContinuous

$newHeaderCellsGJ = $ws.Range("G1:J1")
$newHeaderCellsGJ.Font.Bold = $true
$newHeaderCellsGJ.HorizontalAlignment = -4108  # xlCenter
$newHeaderCellsGJ.VerticalAlignment = -4160    # xlTop
$newHeaderCellsGJ.Borders.LineStyle = 1        # xlContinuous

$lastRow = 13

# --- New leading "type" column for every data row --------------------------
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = "generation"
}

# --- New trailing error-count columns (all zero for this dataset) ---------
$ws.Range("G2:J$lastRow").Value = 0

Write-Host "applied evaluation-by-extraction schema sync"
